# DeltaTime in Layer::OnUpdate() function
# Add two new task rows to the tracker sheet:
#   Row 14: Timer          | (blank subtask) | Done
#   Row 15: App             | FPS; DeltaTime  | Done

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Timer"
$ws.Range("C14").Value = "Done"

$ws.Range("A15").Value = "App"
$ws.Range("B15").Value = "FPS; DeltaTime"
$ws.Range("C15").Value = "Done"

# Leave the selection where the user ended up after entering the data
$ws.Range("B16").Select() | Out-Null
